$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking strings (e.g. "9.30", "0.136")
# are preserved exactly as typed, matching the source inline-string cells.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '66.610.22'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '3.591.83'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '609.17'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('D6').Value = '147.84'
$ws.Range('E6').Value = '  +1.85%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.489'
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('B9').Value = 'Toncoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D9').Value = '8.05'
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.136'
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').Value = '4.201.62'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').Value = '29.93'
$ws.Range('E14').Value = '  -0.85%  '
$ws.Range('D15').Value = '3.589.85'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = '66.664.44'
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = '11.51'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').Value = '6.33'
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').Value = '15.05'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').Value = '428.07'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('D23').Value = '78.84'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').Value = '3.736.49'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  +1.65%  '
$ws.Range('D27').Value = '8.25'
$ws.Range('E27').Value = '  +2.64%  '
$ws.Range('D28').Value = '9.30'
$ws.Range('E28').Value = '  +1.29%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').Value = '3.589.34'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '0.158'
$ws.Range('E32').Value = '  -0.87%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '25.47'
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '1.46'
$ws.Range('E34').Value = '  -1.73%  '
$ws.Range('D35').Value = '7.85'
$ws.Range('E35').Value = '  -0.69%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  -2.99%  '
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('D39').Value = '177.44'
$ws.Range('E39').Value = '  +3.58%  '
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('D41').Value = '5.24'
$ws.Range('E41').Value = '  +0.30%  '
$ws.Range('D42').Value = '0.898'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').Value = '1.91'
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('E44').Value = '  +7.40%  '
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('E46').Value = '  -2.16%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '24.47'
$ws.Range('E47').Value = '  +4.35%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '25.01'
$ws.Range('E48').Value = '  -3.79%  '
$ws.Range('D49').Value = '7.19'
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').Value = '0.235'
$ws.Range('E51').Value = '  -1.66%  '
